$wb = $excel.ActiveWorkbook

# --- Poisson sheet: fill in f(x) and cumulative probability columns (B4:B24, C4:C24) ---
$wsPoisson = $wb.Worksheets.Item("Poisson")
for ($r = 4; $r -le 24; $r++) {
    $wsPoisson.Range("B$r").Formula = "=_xlfn.POISSON.DIST(A$r,`$D`$1,FALSE)"
    $wsPoisson.Range("C$r").Formula = "=_xlfn.POISSON.DIST(A$r,`$D`$1,TRUE)"
}

# --- Normal sheet: add probability formulas ---
$wsNormal = $wb.Worksheets.Item("Normal")
$wsNormal.Range("B5").Formula = "=_xlfn.NORM.DIST(40000,`$B`$1,`$B`$2,TRUE)"
$wsNormal.Range("B6").Formula = "=1-B5"
$wsNormal.Range("B10").Formula = "=_xlfn.NORM.DIST(40000,B1,B2,TRUE)"
$wsNormal.Range("B11").Formula = "=_xlfn.NORM.DIST(30000,B1,B2,TRUE)"
$wsNormal.Range("B12").Formula = "=B10-B11"

# --- Exponential sheet: add probability formulas ---
$wsExponential = $wb.Worksheets.Item("Exponential")
$wsExponential.Range("B3").Formula = "=_xlfn.EXPON.DIST(18,1/B1,TRUE)"
$wsExponential.Range("B4").Formula = "=_xlfn.EXPON.DIST(6,1/B1,TRUE)"
$wsExponential.Range("B5").Formula = "=B3-B4"

# --- Selection / active-tab bookkeeping to match the final saved view state ---
$wsPoisson.Range("C12").Select() | Out-Null
$wsNormal.Range("B13").Select() | Out-Null

$wsExponential.Activate() | Out-Null
$wsExponential.Range("B6").Select() | Out-Null
